$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-16 Thursday" "2025-10-17 Friday"
Replace-Text "839×9=7551" "758×9=6822"
Replace-Text "655×8=5240" "759×9=6831"
Replace-Text "145×7=1015" "222×5=1110"
Replace-Text "515×7=3605" "609×5=3045"
Replace-Text "103×6=618" "476×6=2856"
Replace-Text "142×3=426" "360×6=2160"
Replace-Text "988×6=5928" "832×6=4992"
Replace-Text "125×4=500" "776×9=6984"
Replace-Text "418×9=3762" "177×6=1062"
Replace-Text "459×3=1377" "600×6=3600"
Replace-Text "172×4=688" "526×2=1052"
Replace-Text "773×8=6184" "367×9=3303"
Replace-Text "556×3=1668" "903×6=5418"
Replace-Text "434×7=3038" "562×8=4496"
Replace-Text "485×6=2910" "236×3=708"
Replace-Text "526×9=4734" "264×8=2112"
Replace-Text "233×8=1864" "681×9=6129"
Replace-Text "563×7=3941" "149×8=1192"
Replace-Text "645×5=3225" "497×2=994"
Replace-Text "395×9=3555" "840×8=6720"
Replace-Text "186×4=744" "199×7=1393"
Replace-Text "206×8=1648" "359×7=2513"
Replace-Text "631×4=2524" "305×2=610"
Replace-Text "839×6=5034" "863×9=7767"
Replace-Text "991×7=6937" "486×2=972"
